$wb = $excel.ActiveWorkbook

# Add the new "Tracking" worksheet after the last existing sheet (i.e. at the end)
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$ws = $wb.Worksheets.Add($null, $lastSheet)
$ws.Name = "Tracking"

# Populate cells - order chosen so shared-string insertion order matches the
# target workbook (D3, D2, C2, C3)
$ws.Range("D3").Value = "CSV files substituted for Dimension tables"
$ws.Range("D2").Value = "Commit name"
$ws.Range("C2").Value = "Task"
$ws.Range("C3").Value = "5a"

# Column D width
$ws.Columns.Item(4).ColumnWidth = 45.265625

# Make D3 the active selection on the new sheet (this also becomes the active/tabSelected sheet)
$ws.Range("D3").Select() | Out-Null
